$d = $word.ActiveDocument

# --- 1. Heading "XX. " -> "10. ", with the _GoBack bookmark moved to
#        sit between the number and the ". " (collapsed bookmark range,
#        which forces the run to split exactly like the target markup).
$r = $d.Range(2, 2)
$d.Bookmarks.Add("_GoBack", $r)

$rNum = $d.Range(0, 2)
$rNum.Text = "10"

# (Adding the bookmark above re-targets the single allowed "_GoBack"
#  bookmark, which automatically removes it from its old location
#  further down in the document - no separate deletion step needed.)

# --- 2. Strip the Quick Style flag from the built-in "Normal Table" style.
foreach ($s in $d.Styles) {
    if ($s.NameLocal -eq "Normal Table") {
        $s.QuickStyle = $false
    }
}
